# Fix irrelevant-highlighting issues: adjust a handful of dictionary-entry
# glosses so the "&" separator / "→" normalization lines up correctly.
$d = $word.ActiveDocument

function Fix-Paragraph($index, $find, $replace) {
    $para = $d.Paragraphs.Item($index)
    $range = $para.Range
    $ok = $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 1)
    if (-not $ok) {
        Write-Output "WARNING: replacement failed for paragraph $index"
    }
}

Fix-Paragraph 8 "≈ не творт" "≈ не & творт"
Fix-Paragraph 9 "≈ не сътворт" "≈ не & сътворт"
Fix-Paragraph 17 "въ простѫ → въ + Acc. → въ & простъ" "въ простѫ → въ & простъ"
Fix-Paragraph 37 "творт & ꙁнамен" "# творт ꙁнамен → творт & ꙁнамен"
Fix-Paragraph 47 "въ лѣпотѫ → въ + Acc. → въ & лѣпота" "въ лѣпотѫ → въ & лѣпота"
Fix-Paragraph 48 "по лѣпотѣ → по + Dat. → по & лѣпота" "по лѣпотѣ → по & лѣпота"
Fix-Paragraph 82 "≈ прѧстьнкъ бꙑт" "≈ прѧстьнкъ & бꙑт"
Fix-Paragraph 83 "≈ прьтьнкъ бꙑт" "≈ прьтьнкъ & бꙑт"
Fix-Paragraph 136 "авраамовъ & ѧдь" "авраамова ѧдь → авраамовъ & ѧдь"
Fix-Paragraph 142 "≈ прѧстьнкъ бꙑт" "≈ прѧстьнкъ & бꙑт"
Fix-Paragraph 143 "≈ прьтьнкъ бꙑт" "≈ прьтьнкъ & бꙑт"
